$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date in column C for all data rows (2-391)
#    from 45184 (2023-09-15) to 45186 (2023-09-17).
$ws.Range("C2:C391").Value = 45186

# 2. Add a second argument (the "Beteckning" text, e.g. "A 37268-2021") to the
#    HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2-17 -- the only
#    rows that currently carry those link formulas.
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")
for ($row = 2; $row -le 17; $row++) {
    $beteckning = $ws.Range("A$row").Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula
        if ($formula -and $formula.Length -gt 0) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
            # (keep the comma-space-quote form to mirror Excel's own formula text)
            $cell.Formula = $newFormula
        }
    }
}

# 3. Append two new rows (392 and 393) with new avverkningsanmälningar entries.
$ws.Range("A392").Value = "A 43478-2023"
$ws.Range("B392").Value = 45184
$ws.Range("B392").NumberFormat = "YYYY-MM-DD"
$ws.Range("C392").Value = 45186
$ws.Range("C392").NumberFormat = "YYYY-MM-DD"
$ws.Range("D392").Value = "SKÅNE LÄN"
$ws.Range("E392").Value = "ÖSTRA GÖINGE"
$ws.Range("G392").Value = 1.3
$ws.Range("H392").Value = 0
$ws.Range("I392").Value = 0
$ws.Range("J392").Value = 0
$ws.Range("K392").Value = 0
$ws.Range("L392").Value = 0
$ws.Range("M392").Value = 0
$ws.Range("N392").Value = 0
$ws.Range("O392").Value = 0
$ws.Range("P392").Value = 0
$ws.Range("Q392").Value = 0
$ws.Range("R392").WrapText = $true
$ws.Rows.Item(392).RowHeight = 15

$ws.Range("A393").Value = "A 43479-2023"
$ws.Range("B393").Value = 45184
$ws.Range("B393").NumberFormat = "YYYY-MM-DD"
$ws.Range("C393").Value = 45186
$ws.Range("C393").NumberFormat = "YYYY-MM-DD"
$ws.Range("D393").Value = "SKÅNE LÄN"
$ws.Range("E393").Value = "ÖSTRA GÖINGE"
$ws.Range("G393").Value = 2.2
$ws.Range("H393").Value = 0
$ws.Range("I393").Value = 0
$ws.Range("J393").Value = 0
$ws.Range("K393").Value = 0
$ws.Range("L393").Value = 0
$ws.Range("M393").Value = 0
$ws.Range("N393").Value = 0
$ws.Range("O393").Value = 0
$ws.Range("P393").Value = 0
$ws.Range("Q393").Value = 0
$ws.Range("R393").WrapText = $true

# 4. Row 391 picks up an explicit row height (matches the rest of the sheet)
#    once the sheet is re-saved with the newly appended rows.
$ws.Rows.Item(391).RowHeight = 15
